$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking Price cells to remain text (matches source data which stores
# these as plain strings, e.g. "552.30", not numbers) by temporarily applying a text
# number format, then reverting the style so no extra formatting is left on the cells.
$numericPriceCells = @("D5", "D6", "D8", "D10", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D23", "D25", "D28", "D29", "D32", "D33", "D35", "D36", "D40", "D41", "D42", "D43", "D44", "D45", "D50")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.749.14"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.429.55"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "552.30"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "160.39"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.510"
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("D9").Value = "2.427.64"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  +7.44%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "0.327"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "4.79"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "67.692.98"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "0.0000168"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").Value = "23.04"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "10.30"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").Value = "334.05"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").Value = "6.82"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").Value = "3.77"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("D23").Value = "66.23"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "8.08"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "419.94"
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("E30").Value = "  +2.59%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "160.57"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("D33").Value = "18.93"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "17.79"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "0.103"
$ws.Range("E36").Value = "  -4.98%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "1.07"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "2.00"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").Value = "128.83"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").Value = "0.0710"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "0.478"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("E49").Value = "  -4.94%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "16.58"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0203"
$ws.Range("E51").Value = "  +5.42%  "

foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
